$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 changes from "This is an invalid workbook" to "inventory work book"
$ws.Range("A1").Value = "inventory work book"

# A3, A7, A11 keep displaying "test" (shared-string index churn only, no visible change)
$ws.Range("A3").Value = "test"
$ws.Range("A7").Value = "test"
$ws.Range("A11").Value = "test"

# Selection moves from A16 to A2
$ws.Range("A2").Select()
